$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.015.39'
$ws.Range('E2').Value = '  +1.49%  '
$ws.Range('D3').Value = '3.316.62'
$ws.Range('E3').Value = '  +6.04%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '601.80'
$ws.Range('E5').Value = '  +1.17%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.64'
$ws.Range('E6').Value = '  +5.27%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '3.313.48'
$ws.Range('E8').Value = '  +6.21%  '
$ws.Range('E9').Value = '  +1.48%  '
$ws.Range('E10').Value = '  +3.37%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.54'
$ws.Range('E11').Value = '  +5.75%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.475'
$ws.Range('E12').Value = '  +4.29%  '
$ws.Range('E13').Value = '  +1.67%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.99'
$ws.Range('E14').Value = '  +2.55%  '
$ws.Range('D15').Value = '3.863.41'
$ws.Range('E15').Value = '  +6.11%  '
$ws.Range('E16').Value = '  +0.11%  '
$ws.Range('D17').Value = '3.318.47'
$ws.Range('E17').Value = '  +6.21%  '
$ws.Range('D18').Value = '64.124.53'
$ws.Range('E18').Value = '  +1.64%  '
$ws.Range('E19').Value = '  +3.74%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '482.97'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.33'
$ws.Range('E21').Value = '  +1.33%  '
$ws.Range('E22').Value = '  +6.22%  '
$ws.Range('E23').Value = '  +4.34%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.68'
$ws.Range('E24').Value = '  +5.66%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '85.06'
$ws.Range('E25').Value = '  -1.73%  '
$ws.Range('E26').Value = '  +0.19%  '
$ws.Range('E27').Value = '  +2.37%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.32'
$ws.Range('E28').Value = '  +2.81%  '
$ws.Range('E29').Value = '  -0.11%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.27'
$ws.Range('E30').Value = '  +4.23%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '29.67'
$ws.Range('E31').Value = '  +11.13%  '
$ws.Range('E32').Value = '  +5.73%  '
$ws.Range('E33').Value = '  -1.59%  '
$ws.Range('E34').Value = '  +2.47%  '
$ws.Range('E35').Value = '  +2.43%  '
$ws.Range('E36').Value = '  +3.71%  '
$ws.Range('D37').Value = '0.0₃0762'
$ws.Range('E37').Value = '  +7.38%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '53.35'
$ws.Range('E38').Value = '  +2.58%  '
$ws.Range('E39').Value = '  +4.53%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '435.14'
$ws.Range('E40').Value = '  +2.77%  '
$ws.Range('D41').Value = '3.069.41'
$ws.Range('E41').Value = '  +5.96%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.80'
$ws.Range('E42').Value = '  +3.92%  '
$ws.Range('E43').Value = '  +3.00%  '
$ws.Range('E44').Value = '  -0.69%  '
$ws.Range('E45').Value = '  +2.56%  '
$ws.Range('E46').Value = '  +4.87%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '26.64'
$ws.Range('E47').Value = '  +4.49%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '36.01'
$ws.Range('E48').Value = '  +14.60%  '
$ws.Range('E50').Value = '  +3.03%  '
$ws.Range('E51').Value = '  +1.96%  '
